# Insert a new price-report row before the existing row 182 (1a nueva(o) de
# Pehuenche, Región de La Araucanía) on the "Hortaliza, Feria Lagunitas de
# Puerto Montt - Papa" sheet. All subsequent rows shift down by one, and the
# used range grows from A1:R268 to A1:R269.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 182:268 down to 183:269, opening up a blank row 182.
$ws.Rows("182:182").Insert()

# Populate the newly opened row 182 with the new observation.
$ws.Cells.Item(182, 1).Value = 4
$ws.Cells.Item(182, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(182, 3).Value = "Los Lagos"
$ws.Cells.Item(182, 4).Value = 44510
$ws.Cells.Item(182, 5).Value = 10
$ws.Cells.Item(182, 6).Value = 100114001
$ws.Cells.Item(182, 7).Value = "Papa"
$ws.Cells.Item(182, 8).Value = "Pehuenche"
$ws.Cells.Item(182, 9).Value = "1a nueva(o)"
$ws.Cells.Item(182, 10).Value = 150
$ws.Cells.Item(182, 11).Value = 16000
$ws.Cells.Item(182, 12).Value = 16000
$ws.Cells.Item(182, 13).Value = 16000
$ws.Cells.Item(182, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(182, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(182, 16).Value = 640
$ws.Cells.Item(182, 17).Value = 25
$ws.Cells.Item(182, 18).Value = "Hortaliza"
